$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the RAD Notice Number Error Message timestamps (text values, not dates)
$ws.Range("B2").Value = "Wed Dec 20 12:58:05 EST 2023"
$ws.Range("B4").Value = "Wed Dec 20 12:58:19 EST 2023"
